# Auto-generated edit script applying the Lamia_Profits.xlsx diff
# Updates numeric cell values across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 909935.8
$ws.Range("I6").Value = 1111821.5
$ws.Range("K6").Value = 3335464.5
$ws.Range("M6").Value = -3335352.5
$ws.Range("H12").Value = 183.16667
$ws.Range("I12").Value = 155.55556
$ws.Range("J12").Value = 266
$ws.Range("K12").Value = 155.55556
$ws.Range("L12").Value = 266
$ws.Range("M12").Value = 14.44443999999999
$ws.Range("N12").Value = -606
$ws.Range("H15").Value = 1985444.8
$ws.Range("I15").Value = 1985444.8
$ws.Range("K15").Value = 5956334.4
$ws.Range("M15").Value = -5956165.4
$ws.Range("H29").Value = 3704
$ws.Range("J29").Value = 7953
$ws.Range("L29").Value = 23859
$ws.Range("N29").Value = -24421
$ws.Range("H39").Value = 295
$ws.Range("I39").Value = 211.63637
$ws.Range("K39").Value = 634.9091100000001
$ws.Range("M39").Value = -338.9091100000001
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H47").Value = 22350
$ws.Range("J47").Value = 35000
$ws.Range("L47").Value = 35000
$ws.Range("N47").Value = -36944
$ws.Range("H48").Value = 8832.134
$ws.Range("J48").Value = 10623
$ws.Range("L48").Value = 31869
$ws.Range("N48").Value = -32453
$ws.Range("H51").Value = 8941.333000000001
$ws.Range("J51").Value = 9762.5
$ws.Range("L51").Value = 9762.5
$ws.Range("N51").Value = -10730.5
$ws.Range("H53").Value = 2092.6155
$ws.Range("I53").Value = 326
$ws.Range("J53").Value = 2877.7778
$ws.Range("K53").Value = 326
$ws.Range("L53").Value = 2877.7778
$ws.Range("M53").Value = 311
$ws.Range("N53").Value = -4151.7778
$ws.Range("H54").Value = 5300
$ws.Range("I54").Value = 5300
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 5300
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -4814
$ws.Range("N54").ClearContents()
$ws.Range("H55").Value = 898.8946999999999
$ws.Range("J55").Value = 1562.875
$ws.Range("L55").Value = 1562.875
$ws.Range("N55").Value = -1990.875
$ws.Range("H56").Value = 8832.134
$ws.Range("J56").Value = 10623
$ws.Range("L56").Value = 31869
$ws.Range("N56").Value = -32937
$ws.Range("H58").Value = 871.6667
$ws.Range("I58").Value = 83.5
$ws.Range("J58").Value = 2448
$ws.Range("K58").Value = 250.5
$ws.Range("L58").Value = 7344
$ws.Range("M58").Value = -100.5
$ws.Range("N58").Value = -7644
$ws.Range("H61").Value = 1922.1111
$ws.Range("I61").Value = 1922.1111
$ws.Range("K61").Value = 5766.3333
$ws.Range("M61").Value = -5594.3333
$ws.Range("H64").Value = 9146.923000000001
$ws.Range("J64").Value = 9075.833000000001
$ws.Range("L64").Value = 9075.833000000001
$ws.Range("N64").Value = -9571.833000000001
$ws.Range("H67").Value = 9146.923000000001
$ws.Range("J67").Value = 9075.833000000001
$ws.Range("L67").Value = 9075.833000000001
$ws.Range("N67").Value = -10791.833
$ws.Range("H74").Value = 7125.1904
$ws.Range("I74").Value = 5164.0625
$ws.Range("K74").Value = 5164.0625
$ws.Range("M74").Value = -4228.0625
$ws.Range("H77").Value = 7125.1904
$ws.Range("I77").Value = 5164.0625
$ws.Range("K77").Value = 25820.3125
$ws.Range("M77").Value = -21140.3125
$ws.Range("H80").Value = 2791.5625
$ws.Range("I80").Value = 587.1667
$ws.Range("K80").Value = 1761.5001
$ws.Range("M80").Value = -763.5001
$ws.Range("H83").Value = 2791.5625
$ws.Range("I83").Value = 587.1667
$ws.Range("K83").Value = 5284.5003
$ws.Range("M83").Value = -292.5002999999997
$ws.Range("H88").Value = 5455.9287
$ws.Range("J88").Value = 5567.923
$ws.Range("L88").Value = 5567.923
$ws.Range("N88").Value = -6379.923
$ws.Range("H91").Value = 5455.9287
$ws.Range("J91").Value = 5567.923
$ws.Range("L91").Value = 5567.923
$ws.Range("N91").Value = -8375.922999999999
$ws.Range("H100").Value = 106374.6
$ws.Range("I100").Value = 147599.86
$ws.Range("K100").Value = 147599.86
$ws.Range("M100").Value = -147058.86
$ws.Range("H107").Value = 1829.8462
$ws.Range("I107").Value = 1748.5
$ws.Range("K107").Value = 1748.5
$ws.Range("M107").Value = 171.5
$ws.Range("H111").Value = 823
$ws.Range("I111").Value = 655.6667
$ws.Range("J111").Value = 1023.8
$ws.Range("K111").Value = 1967.0001
$ws.Range("L111").Value = 3071.4
$ws.Range("M111").Value = 1099.9999
$ws.Range("N111").Value = -9205.4
$ws.Range("H112").Value = 2274.8667
$ws.Range("J112").Value = 2561.3635
$ws.Range("L112").Value = 7684.0905
$ws.Range("N112").Value = -9900.0905
$ws.Range("H118").Value = 1273.6666
$ws.Range("I118").Value = 514.25
$ws.Range("J118").Value = 1881.2
$ws.Range("K118").Value = 1542.75
$ws.Range("L118").Value = 5643.6
$ws.Range("M118").Value = 114.25
$ws.Range("N118").Value = -8957.6
$ws.Range("H125").Value = 1380.6364
$ws.Range("I125").Value = 1466
$ws.Range("J125").Value = 1361.6666
$ws.Range("K125").Value = 13194
$ws.Range("L125").Value = 12254.9994
$ws.Range("M125").Value = -10734
$ws.Range("N125").Value = -17174.9994
$ws.Range("H127").Value = 1202.1428
$ws.Range("I127").Value = 483.5
$ws.Range("J127").Value = 2998.75
$ws.Range("K127").Value = 1450.5
$ws.Range("L127").Value = 8996.25
$ws.Range("M127").Value = 3509.5
$ws.Range("N127").Value = -18916.25
$ws.Range("H132").Value = 1970.6666
$ws.Range("I132").Value = 1963.3334
$ws.Range("J132").Value = 2029.3334
$ws.Range("K132").Value = 5890.0002
$ws.Range("L132").Value = 6088.0002
$ws.Range("M132").Value = -3360.0002
$ws.Range("N132").Value = -11148.0002
$ws.Range("H137").Value = 3273.4324
$ws.Range("I137").Value = 2464.625
$ws.Range("J137").Value = 3889.6667
$ws.Range("K137").Value = 7393.875
$ws.Range("L137").Value = 11669.0001
$ws.Range("M137").Value = -4843.875
$ws.Range("N137").Value = -16769.0001
$ws.Range("H138").Value = 3262.2
$ws.Range("I138").Value = 2856.2307
$ws.Range("J138").Value = 3387.8572
$ws.Range("K138").Value = 8568.6921
$ws.Range("L138").Value = 10163.5716
$ws.Range("M138").Value = -3428.6921
$ws.Range("N138").Value = -20443.5716
$ws.Range("H141").Value = 822.5625
$ws.Range("J141").Value = 505
$ws.Range("L141").Value = 1515
$ws.Range("N141").Value = -11875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1771.0128
$ws.Range("I32").Value = 1636.3334
$ws.Range("J32").Value = 5138
$ws.Range("K32").Value = 1636.3334
$ws.Range("L32").Value = 5138
$ws.Range("M32").Value = -1349.3334
$ws.Range("N32").Value = -5712
$ws.Range("H96").Value = 47337.8
$ws.Range("J96").Value = 47337.8
$ws.Range("L96").Value = 47337.8
$ws.Range("N96").Value = -52829.8
$ws.Range("H122").Value = 83336100
$ws.Range("I122").Value = 3015.182
$ws.Range("K122").Value = 9045.545999999998
$ws.Range("M122").Value = -6595.545999999998
$ws.Range("H132").Value = 8205.966
$ws.Range("I132").Value = 4109.294
$ws.Range("K132").Value = 12327.882
$ws.Range("M132").Value = -9797.882

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 5011.3335
$ws.Range("I24").Value = 2525
$ws.Range("J24").Value = 9984
$ws.Range("K24").Value = 2525
$ws.Range("L24").Value = 9984
$ws.Range("M24").Value = -2290
$ws.Range("N24").Value = -10454
$ws.Range("H95").Value = 17280.75
$ws.Range("J95").Value = 17280.75
$ws.Range("L95").Value = 17280.75
$ws.Range("N95").Value = -22772.75
$ws.Range("H97").Value = 24241.428
$ws.Range("J97").Value = 55179
$ws.Range("L97").Value = 55179
$ws.Range("N97").Value = -57161
$ws.Range("H99").Value = 2306.6667
$ws.Range("I99").Value = 2306.6667
$ws.Range("K99").Value = 2306.6667
$ws.Range("M99").Value = -808.6667000000002
$ws.Range("H105").Value = 15460.85
$ws.Range("I105").Value = 14637.9375
$ws.Range("K105").Value = 14637.9375
$ws.Range("M105").Value = -12890.9375
$ws.Range("H107").Value = 1613.45
$ws.Range("I107").Value = 1329.9474
$ws.Range("K107").Value = 1329.9474
$ws.Range("M107").Value = 590.0526
$ws.Range("H138").Value = 65490.57
$ws.Range("J138").Value = 65490.57
$ws.Range("L138").Value = 65490.57
$ws.Range("N138").Value = -75770.57000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 72500.92999999999
$ws.Range("J31").Value = 134690.88
$ws.Range("L31").Value = 134690.88
$ws.Range("N31").Value = -135280.88
$ws.Range("H34").Value = 72500.92999999999
$ws.Range("J34").Value = 134690.88
$ws.Range("L34").Value = 134690.88
$ws.Range("N34").Value = -135094.88
$ws.Range("H99").Value = 3212.3076
$ws.Range("I99").Value = 3118.1428
$ws.Range("K99").Value = 3118.1428
$ws.Range("M99").Value = -1620.1428
$ws.Range("H126").Value = 3212.3076
$ws.Range("I126").Value = 3118.1428
$ws.Range("K126").Value = 9354.428400000001
$ws.Range("M126").Value = -6884.428400000001
$ws.Range("H132").Value = 2614.681
$ws.Range("I132").Value = 2035.3334
$ws.Range("J132").Value = 5439
$ws.Range("K132").Value = 6106.0002
$ws.Range("L132").Value = 16317
$ws.Range("M132").Value = -3576.0002
$ws.Range("N132").Value = -21377
$ws.Range("H134").Value = 3777.805
$ws.Range("I134").Value = 2413.9688
$ws.Range("K134").Value = 7241.9064
$ws.Range("M134").Value = -4706.9064

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8196845
$ws.Range("I4").Value = 6902679
$ws.Range("K4").Value = 20708037
$ws.Range("M4").Value = -20707925
$ws.Range("H33").Value = 909147.4
$ws.Range("I33").Value = 1111154.5
$ws.Range("J33").Value = 115
$ws.Range("K33").Value = 6666927
$ws.Range("L33").Value = 690
$ws.Range("M33").Value = -6666644
$ws.Range("N33").Value = -1256
$ws.Range("H41").Value = 250
$ws.Range("J41").Value = 300
$ws.Range("L41").Value = 900
$ws.Range("N41").Value = -1576
$ws.Range("H88").Value = 11617.5
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 11617.5
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H106").Value = 13026.667
$ws.Range("J106").Value = 19027
$ws.Range("L106").Value = 57081
$ws.Range("N106").Value = -58973
$ws.Range("H113").Value = 1610
$ws.Range("I113").Value = 1849
$ws.Range("J113").Value = 1566.5454
$ws.Range("K113").Value = 5547
$ws.Range("L113").Value = 4699.6362
$ws.Range("M113").Value = -3377
$ws.Range("N113").Value = -9039.636200000001
$ws.Range("H132").Value = 3640.3
$ws.Range("I132").Value = 2748.3333
$ws.Range("J132").Value = 4532.2666
$ws.Range("K132").Value = 24734.9997
$ws.Range("L132").Value = 40790.3994
$ws.Range("M132").Value = -22204.9997
$ws.Range("N132").Value = -45850.3994
$ws.Range("H140").Value = 3749.5
$ws.Range("I140").Value = 1625
$ws.Range("K140").Value = 4875
$ws.Range("M140").Value = 305

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14155.214
$ws.Range("I70").Value = 10307
$ws.Range("J70").Value = 18003.428
$ws.Range("K70").Value = 10307
$ws.Range("L70").Value = 18003.428
$ws.Range("M70").Value = -10037
$ws.Range("N70").Value = -18543.428
$ws.Range("H73").Value = 14155.214
$ws.Range("I73").Value = 10307
$ws.Range("J73").Value = 18003.428
$ws.Range("K73").Value = 10307
$ws.Range("L73").Value = 18003.428
$ws.Range("M73").Value = -9371
$ws.Range("N73").Value = -19875.428
$ws.Range("H92").Value = 16720
$ws.Range("J92").Value = 16960
$ws.Range("L92").Value = 16960
$ws.Range("N92").Value = -20704
$ws.Range("H102").Value = 2124.8696
$ws.Range("I102").Value = 2124.8696
$ws.Range("K102").Value = 2124.8696
$ws.Range("M102").Value = -502.8696
$ws.Range("H107").Value = 1071.3334
$ws.Range("I107").Value = 295.1111
$ws.Range("J107").Value = 3400
$ws.Range("K107").Value = 295.1111
$ws.Range("L107").Value = 3400
$ws.Range("M107").Value = 1624.8889
$ws.Range("N107").Value = -7240
$ws.Range("H122").Value = 7358.579
$ws.Range("I122").Value = 8172.643
$ws.Range("K122").Value = 24517.929
$ws.Range("M122").Value = -22067.929
$ws.Range("H132").Value = 348283.38
$ws.Range("I132").Value = 412589.22
$ws.Range("K132").Value = 1237767.66
$ws.Range("M132").Value = -1235237.66

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5894.8945
$ws.Range("I7").Value = 3671.9375
$ws.Range("J7").Value = 17750.666
$ws.Range("K7").Value = 3671.9375
$ws.Range("L7").Value = 17750.666
$ws.Range("M7").Value = -3559.9375
$ws.Range("N7").Value = -17974.666
$ws.Range("H46").Value = 3028.2666
$ws.Range("J46").Value = 3101.7144
$ws.Range("L46").Value = 3101.7144
$ws.Range("N46").Value = -3477.7144
$ws.Range("H68").Value = 94582
$ws.Range("I68").Value = 113133.336
$ws.Range("K68").Value = 113133.336
$ws.Range("M68").Value = -112384.336
$ws.Range("H71").Value = 94582
$ws.Range("I71").Value = 113133.336
$ws.Range("K71").Value = 565666.6799999999
$ws.Range("M71").Value = -561922.6799999999
$ws.Range("H82").Value = 5077.7144
$ws.Range("I82").Value = 4415.5835
$ws.Range("J82").Value = 5574.3125
$ws.Range("K82").Value = 4415.5835
$ws.Range("L82").Value = 5574.3125
$ws.Range("M82").Value = -4054.5835
$ws.Range("N82").Value = -6296.3125
$ws.Range("H85").Value = 5077.7144
$ws.Range("I85").Value = 4415.5835
$ws.Range("J85").Value = 5574.3125
$ws.Range("K85").Value = 4415.5835
$ws.Range("L85").Value = 5574.3125
$ws.Range("M85").Value = -3167.5835
$ws.Range("N85").Value = -8070.3125
$ws.Range("H100").Value = 12250.823
$ws.Range("I100").Value = 10071.048
$ws.Range("J100").Value = 15772
$ws.Range("K100").Value = 10071.048
$ws.Range("L100").Value = 15772
$ws.Range("M100").Value = -9530.048000000001
$ws.Range("N100").Value = -16854
$ws.Range("H126").Value = 5894.8945
$ws.Range("I126").Value = 3671.9375
$ws.Range("J126").Value = 17750.666
$ws.Range("K126").Value = 11015.8125
$ws.Range("L126").Value = 53251.99800000001
$ws.Range("M126").Value = -8545.8125
$ws.Range("N126").Value = -58191.99800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 33999.285
$ws.Range("I70").Value = 32832.832
$ws.Range("J70").Value = 40998
$ws.Range("K70").Value = 32832.832
$ws.Range("L70").Value = 40998
$ws.Range("M70").Value = -32517.832
$ws.Range("N70").Value = -41628
$ws.Range("H73").Value = 33999.285
$ws.Range("I73").Value = 32832.832
$ws.Range("J73").Value = 40998
$ws.Range("K73").Value = 32832.832
$ws.Range("L73").Value = 40998
$ws.Range("M73").Value = -31740.832
$ws.Range("N73").Value = -43182
$ws.Range("H81").Value = 4152.769
$ws.Range("I81").Value = 2915.3333
$ws.Range("K81").Value = 5830.6666
$ws.Range("M81").Value = -4769.6666
$ws.Range("H84").Value = 4152.769
$ws.Range("I84").Value = 2915.3333
$ws.Range("K84").Value = 29153.333
$ws.Range("M84").Value = -23849.333
$ws.Range("H95").Value = 30344
$ws.Range("J95").Value = 30344
$ws.Range("L95").Value = 30344
$ws.Range("N95").Value = -35836
$ws.Range("H97").Value = 8500.5
$ws.Range("J97").Value = 8500.5
$ws.Range("L97").Value = 8500.5
$ws.Range("N97").Value = -10482.5
$ws.Range("H100").Value = 954.8
$ws.Range("I100").Value = 742.73334
$ws.Range("J100").Value = 1591
$ws.Range("K100").Value = 1485.46668
$ws.Range("L100").Value = 3182
$ws.Range("M100").Value = -944.46668
$ws.Range("N100").Value = -4264
$ws.Range("H107").Value = 1199.2858
$ws.Range("I107").Value = 1483.6666
$ws.Range("K107").Value = 4450.9998
$ws.Range("M107").Value = -2530.9998
$ws.Range("H108").Value = 77473.5
$ws.Range("J108").Value = 77473.5
$ws.Range("L108").Value = 77473.5
$ws.Range("N108").Value = -85153.5
$ws.Range("H113").Value = 808.44446
$ws.Range("I113").Value = 394.5
$ws.Range("J113").Value = 1139.6
$ws.Range("K113").Value = 1183.5
$ws.Range("L113").Value = 3418.8
$ws.Range("M113").Value = 986.5
$ws.Range("N113").Value = -7758.799999999999
$ws.Range("H126").Value = 2696
$ws.Range("J126").Value = 3242
$ws.Range("L126").Value = 9726
$ws.Range("N126").Value = -14666

Write-Host "Applied 431 cell updates across 8 sheets"
